$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New transaction rows appended to the log (rows 86-100).
# Columns A-C are text (date/time string, phone number, name); set as Text
# format first so Excel doesn't auto-convert the numeric-looking strings,
# then restore the default style so no extra formatting is left behind.
$newRange = $ws.Range("A86:C100")
$newRange.NumberFormat = "@"

$ws.Range("A86").Value = '2026-02-14 12:28:29'
$ws.Range("B86").Value = '237671823369'
$ws.Range("C86").Value = 'MFS ENTREE COLLEGE MALANGUE'
$ws.Range("D86").Value = 157448

$ws.Range("A87").Value = '2026-02-14 12:08:46'
$ws.Range("B87").Value = '237672128028'
$ws.Range("C87").Value = 'CAROLINE WAKO DJAMNOU'
$ws.Range("D87").Value = 65097

$ws.Range("A88").Value = '2026-02-14 16:56:53'
$ws.Range("B88").Value = '237672277367'
$ws.Range("C88").Value = 'TOP MOBIL KM5 LTDLA_POLAS_BTQ_KM5'
$ws.Range("D88").Value = 2996

$ws.Range("A89").Value = '2026-02-14 16:47:59'
$ws.Range("B89").Value = '237674853971'
$ws.Range("C89").Value = 'NJOSSEU TCHOUNZOU TOP MOBILE'
$ws.Range("D89").Value = 4306

$ws.Range("A90").Value = '2026-02-14 17:33:03'
$ws.Range("B90").Value = '237674884705'
$ws.Range("C90").Value = 'BAH AMADOU MOUNTAGHA ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Range("D90").Value = 9333

$ws.Range("A91").Value = '2026-02-14 17:10:24'
$ws.Range("B91").Value = '237675779272'
$ws.Range("C91").Value = 'RODES NGWEM KEMAYOU'
$ws.Range("D91").Value = 30585

$ws.Range("A92").Value = '2026-02-13 15:13:34'
$ws.Range("B92").Value = '237677304210'
$ws.Range("C92").Value = 'CARINE OROCK'
$ws.Range("D92").Value = 235659

$ws.Range("A93").Value = '2026-02-06 09:50:11'
$ws.Range("B93").Value = '237678267353'
$ws.Range("C93").Value = 'LA NEGRESSE SARL EMBOLA BELTUS MBU'
$ws.Range("D93").Value = 0

$ws.Range("A94").Value = '2026-02-14 14:35:06'
$ws.Range("B94").Value = '237678370615'
$ws.Range("C94").Value = 'ESSEN ONGOLONG BERTHE HORTENSE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Range("D94").Value = 510310

$ws.Range("A95").Value = '2026-02-14 10:53:33'
$ws.Range("B95").Value = '237678836319'
$ws.Range("C95").Value = 'KAMDOM DOMINIQUE STEPHANIE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Range("D95").Value = 137271

$ws.Range("A96").Value = '2026-02-14 13:39:03'
$ws.Range("B96").Value = '237678922502'
$ws.Range("C96").Value = 'NWOAGA TCHAMDJOU EPSE KAMSEU EMILINE ETS LE CONTENT'
$ws.Range("D96").Value = 730908

$ws.Range("A97").Value = '2026-02-14 08:39:32'
$ws.Range("B97").Value = '237679884264'
$ws.Range("C97").Value = 'GABRIEL MONKAM TCHOUPE'
$ws.Range("D97").Value = 365895

$ws.Range("A98").Value = '2026-02-14 11:35:46'
$ws.Range("B98").Value = '237681019523'
$ws.Range("C98").Value = 'ETS MOULAY RIPERT AND COMPANY'
$ws.Range("D98").Value = 245762

$ws.Range("A99").Value = '2026-02-14 17:09:47'
$ws.Range("B99").Value = '237681125655'
$ws.Range("C99").Value = 'EMENGUE PICHOU ROMEO KAMILAH CONNECTION GROUP'
$ws.Range("D99").Value = 33916

$ws.Range("A100").Value = '2026-02-14 09:46:40'
$ws.Range("B100").Value = '237681240793'
$ws.Range("C100").Value = 'MBANE EMILIE FRANCOISE ETS MOBILE FINANCIAL SERVICES MFS'
$ws.Range("D100").Value = 2067

$newRange.Style = "Normal"

Write-Host "Dimension / UsedRange rows:" $ws.UsedRange.Rows.Count
